# TaskList.xlsx — "UI" sheet tidy-up + content update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UI")

# ---------------------------------------------------------------------------
# 1. Content changes
# ---------------------------------------------------------------------------
# Row 4 ("Number of photos retrieved ...") is now marked Pass in the Testing column.
$ws.Range("H5").Value = "Pass"

# Clarify the bug description: "no pop up" -> "no pop up/error"
$ws.Range("B14").Value = 'When the directory/date folder entered is not available, there is no pop up/error saying "no such directory exists"'

# ---------------------------------------------------------------------------
# 2. Whole-table formatting clean-up: a uniform thin border box around every
#    used cell (A1:H14), replacing the old mix of thin/medium borders.
# ---------------------------------------------------------------------------
$full = $ws.Range("A1:H14")
$full.Borders.LineStyle = -4142   # xlLineStyleNone - start every cell from a clean slate
$full.Borders.LineStyle = 1       # xlContinuous
$full.Borders.Weight = 2          # xlThin

# Header row: no more thick bottom divider / extra row height, just bold text.
$ws.Rows("1:1").RowHeight = 14.5
$hdr = $ws.Range("A1:H1")
$hdr.Font.Bold = $true
$hdr.WrapText = $false
$hdr.HorizontalAlignment = -4131  # xlGeneral

# Column A (row numbers): centered.
$ws.Range("A2:A14").HorizontalAlignment = -4108  # xlCenter

# Column B (issue descriptions): wrap text.
$ws.Range("B2:B14").WrapText = $true

# ---------------------------------------------------------------------------
# 3. Selection left on H5 (the newly-updated cell) when the file was saved.
# ---------------------------------------------------------------------------
$ws.Range("H5").Select()
